$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 559.2727
$ws.Range("I4").Value = 433.77777
$ws.Range("K4").Value = 433.77777
$ws.Range("M4").Value = -319.77777

$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws.Range("H12").Value = 223
$ws.Range("I12").Value = 179.16667
$ws.Range("J12").Value = 275.6
$ws.Range("K12").Value = 179.16667
$ws.Range("L12").Value = 275.6
$ws.Range("M12").Value = -9.166670000000011
$ws.Range("N12").Value = -615.6

$ws.Range("H63").Value = 50000
$ws.Range("J63").Value = 50000
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51248

$ws.Range("H66").Value = 50000
$ws.Range("J66").Value = 50000
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -156240

$ws.Range("H132").Value = 6614.5312
$ws.Range("J132").Value = 3722.375
$ws.Range("L132").Value = 11167.125
$ws.Range("N132").Value = -16227.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2785.7646
$ws.Range("I32").Value = 2785.7646
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2785.7646
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2498.7646
$ws.Range("N32").ClearContents()

$ws.Range("H45").Value = 25502.723
$ws.Range("I45").Value = 33473.54
$ws.Range("K45").Value = 33473.54
$ws.Range("M45").Value = -33096.54

$ws.Range("H61").Value = 4423.533
$ws.Range("I61").Value = 2213.2273
$ws.Range("K61").Value = 2213.2273
$ws.Range("M61").Value = -2001.2273

$ws.Range("H74").Value = 160333.58
$ws.Range("I74").Value = 265554.94
$ws.Range("K74").Value = 265554.94
$ws.Range("M74").Value = -264680.94

$ws.Range("H77").Value = 160333.58
$ws.Range("I77").Value = 265554.94
$ws.Range("K77").Value = 1327774.7
$ws.Range("M77").Value = -1323406.7

$ws.Range("H88").Value = 3192.7273
$ws.Range("I88").Value = 812.75
$ws.Range("J88").Value = 4552.7144
$ws.Range("K88").Value = 812.75
$ws.Range("L88").Value = 4552.7144
$ws.Range("M88").Value = -406.75
$ws.Range("N88").Value = -5364.7144

$ws.Range("H91").Value = 3192.7273
$ws.Range("I91").Value = 812.75
$ws.Range("J91").Value = 4552.7144
$ws.Range("K91").Value = 812.75
$ws.Range("L91").Value = 4552.7144
$ws.Range("M91").Value = 591.25
$ws.Range("N91").Value = -7360.7144

$ws.Range("H110").Value = 1195.2069
$ws.Range("J110").Value = 4560.6
$ws.Range("L110").Value = 4560.6
$ws.Range("N110").Value = -8650.6

$ws.Range("H122").Value = 3237.4443
$ws.Range("I122").Value = 3525.7334
$ws.Range("J122").Value = 1796
$ws.Range("K122").Value = 10577.2002
$ws.Range("L122").Value = 5388
$ws.Range("M122").Value = -8127.200199999999
$ws.Range("N122").Value = -10288

$ws.Range("H136").Value = 4423.533
$ws.Range("I136").Value = 2213.2273
$ws.Range("K136").Value = 6639.6819
$ws.Range("M136").Value = -4089.6819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 971
$ws.Range("J64").Value = 941
$ws.Range("L64").Value = 941
$ws.Range("N64").Value = -1391

$ws.Range("H67").Value = 971
$ws.Range("J67").Value = 941
$ws.Range("L67").Value = 941
$ws.Range("N67").Value = -2501

$ws.Range("H107").Value = 1597.4615
$ws.Range("I107").Value = 1439.0625
$ws.Range("K107").Value = 1439.0625
$ws.Range("M107").Value = 480.9375

$ws.Range("H134").Value = 2836.3584
$ws.Range("I134").Value = 2533
$ws.Range("J134").Value = 4140.8
$ws.Range("K134").Value = 7599
$ws.Range("L134").Value = 12422.4
$ws.Range("M134").Value = -5064
$ws.Range("N134").Value = -17492.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1327.3
$ws.Range("I22").Value = 1434.125
$ws.Range("K22").Value = 1434.125
$ws.Range("M22").Value = -1084.125

$ws.Range("H31").Value = 3135.76
$ws.Range("I31").Value = 2227.1538
$ws.Range("J31").Value = 6357.1816
$ws.Range("K31").Value = 2227.1538
$ws.Range("L31").Value = 6357.1816
$ws.Range("M31").Value = -1932.1538
$ws.Range("N31").Value = -6947.1816

$ws.Range("H34").Value = 3135.76
$ws.Range("I34").Value = 2227.1538
$ws.Range("J34").Value = 6357.1816
$ws.Range("K34").Value = 2227.1538
$ws.Range("L34").Value = 6357.1816
$ws.Range("M34").Value = -2025.1538
$ws.Range("N34").Value = -6761.1816

$ws.Range("H58").Value = 4124.25
$ws.Range("I58").Value = 3874.25
$ws.Range("K58").Value = 3874.25
$ws.Range("M58").Value = -3671.25

$ws.Range("H94").Value = 1836
$ws.Range("I94").Value = 1427.3334
$ws.Range("K94").Value = 1427.3334
$ws.Range("M94").Value = -976.3334

$ws.Range("H122").Value = 2786.7856
$ws.Range("I122").Value = 2391.762
$ws.Range("K122").Value = 7175.286
$ws.Range("M122").Value = -4725.286

$ws.Range("H132").Value = 33336432
$ws.Range("I132").Value = 55556332
$ws.Range("J132").Value = 6582.5
$ws.Range("K132").Value = 166668996
$ws.Range("L132").Value = 19747.5
$ws.Range("M132").Value = -166666466
$ws.Range("N132").Value = -24807.5

$ws.Range("H136").Value = 4124.25
$ws.Range("I136").Value = 3874.25
$ws.Range("K136").Value = 11622.75
$ws.Range("M136").Value = -9072.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 111117200
$ws.Range("J22").Value = 7257
$ws.Range("L22").Value = 21771
$ws.Range("N22").Value = -22109

$ws.Range("H27").Value = 111117200
$ws.Range("J27").Value = 7257
$ws.Range("L27").Value = 21771
$ws.Range("N27").Value = -21975

$ws.Range("H62").Value = 9767.5
$ws.Range("J62").Value = 9767.5
$ws.Range("L62").Value = 29302.5
$ws.Range("N62").Value = -30674.5

$ws.Range("H65").Value = 9767.5
$ws.Range("J65").Value = 9767.5
$ws.Range("L65").Value = 87907.5
$ws.Range("N65").Value = -94771.5

$ws.Range("H81").Value = 4229.5
$ws.Range("J81").Value = 5127.857
$ws.Range("L81").Value = 15383.571
$ws.Range("N81").Value = -17629.571

$ws.Range("H84").Value = 4229.5
$ws.Range("J84").Value = 5127.857
$ws.Range("L84").Value = 46150.713
$ws.Range("N84").Value = -57382.713

$ws.Range("H107").Value = 521.75
$ws.Range("J107").Value = 438.8
$ws.Range("L107").Value = 1316.4
$ws.Range("N107").Value = -5156.4

$ws.Range("H131").Value = 7366.077
$ws.Range("I131").Value = 15688.7
$ws.Range("J131").Value = 2164.4375
$ws.Range("K131").Value = 47066.10000000001
$ws.Range("L131").Value = 6493.3125
$ws.Range("M131").Value = -42026.10000000001
$ws.Range("N131").Value = -16573.3125

$ws.Range("H140").Value = 10810.108
$ws.Range("I140").Value = 8005.2607
$ws.Range("K140").Value = 24015.7821
$ws.Range("M140").Value = -18835.7821

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3762.3447
$ws.Range("I122").Value = 2597.0557
$ws.Range("K122").Value = 7791.1671
$ws.Range("M122").Value = -5341.1671

$ws.Range("H126").Value = 8308.182000000001
$ws.Range("I126").Value = 5082.1665
$ws.Range("J126").Value = 12179.4
$ws.Range("K126").Value = 15246.4995
$ws.Range("L126").Value = 36538.2
$ws.Range("M126").Value = -12776.4995
$ws.Range("N126").Value = -41478.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 691.0909
$ws.Range("I16").Value = 693.0476
$ws.Range("K16").Value = 693.0476
$ws.Range("M16").Value = -523.0476

$ws.Range("H22").Value = 535.46155
$ws.Range("I22").Value = 416
$ws.Range("J22").Value = 610.125
$ws.Range("K22").Value = 416
$ws.Range("L22").Value = 610.125
$ws.Range("M22").Value = -121
$ws.Range("N22").Value = -1200.125

$ws.Range("H27").Value = 535.46155
$ws.Range("I27").Value = 416
$ws.Range("J27").Value = 610.125
$ws.Range("K27").Value = 416
$ws.Range("L27").Value = 610.125
$ws.Range("M27").Value = -309
$ws.Range("N27").Value = -824.125

$ws.Range("H40").Value = 5027.4116
$ws.Range("I40").Value = 5128.154
$ws.Range("J40").Value = 4700
$ws.Range("K40").Value = 5128.154
$ws.Range("L40").Value = 4700
$ws.Range("M40").Value = -4992.154
$ws.Range("N40").Value = -4972

$ws.Range("H46").Value = 2315
$ws.Range("J46").Value = 1549
$ws.Range("L46").Value = 1549
$ws.Range("N46").Value = -1925

$ws.Range("H55").Value = 913.46155
$ws.Range("I55").Value = 496.42856
$ws.Range("K55").Value = 496.42856
$ws.Range("M55").Value = -323.42856
